# Scheduled market-data refresh for the Leve-profit workbook.
#
# Columns on every sheet (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR):
#   H currentAveragePrice    I currentAveragePriceNQ   J currentAveragePriceHQ
#   K LevePriceNQ            L LevePriceHQ             M LeveProfitNQ   N LeveProfitHQ
#
# This run refreshes H:N for the rows whose underlying market prices moved since
# the last pull. A few rows previously showed a stale 0/placeholder profit figure
# (no real NQ listings) -- those profit cells are cleared instead of recomputed,
# and one row that had no NQ-profit value before now gets one.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(6, 8).Value = 290.9091   # H6: 295.5 -> 290.9091
$ws.Cells.Item(6, 10).Value = 260   # J6: 265 -> 260
$ws.Cells.Item(6, 12).Value = 780   # L6: 795 -> 780
$ws.Cells.Item(6, 14).Value = -1004   # N6: -1019 -> -1004

$ws.Cells.Item(40, 8).Value = 6549.6875   # H40: 6899.7334 -> 6549.6875
$ws.Cells.Item(40, 10).Value = 11258.8   # J40: 13748.75 -> 11258.8
$ws.Cells.Item(40, 12).Value = 11258.8   # L40: 13748.75 -> 11258.8
$ws.Cells.Item(40, 14).Value = -11608.8   # N40: -14098.75 -> -11608.8

$ws.Cells.Item(100, 8).Value = 8229.954   # H100: 8203.174000000001 -> 8229.954
$ws.Cells.Item(100, 9).Value = 1706.4   # I100: 1742.6666 -> 1706.4
$ws.Cells.Item(100, 10).Value = 13666.25   # J100: 12356.357 -> 13666.25
$ws.Cells.Item(100, 11).Value = 1706.4   # K100: 1742.6666 -> 1706.4
$ws.Cells.Item(100, 12).Value = 13666.25   # L100: 12356.357 -> 13666.25
$ws.Cells.Item(100, 13).Value = -1165.4   # M100: -1201.6666 -> -1165.4
$ws.Cells.Item(100, 14).Value = -14748.25   # N100: -13438.357 -> -14748.25

$ws.Cells.Item(133, 8).Value = 39464.39   # H133: 39381.906 -> 39464.39
$ws.Cells.Item(133, 10).Value = 39464.39   # J133: 39381.906 -> 39464.39
$ws.Cells.Item(133, 12).Value = 39464.39   # L133: 39381.906 -> 39464.39
$ws.Cells.Item(133, 14).Value = -49584.39   # N133: -49501.906 -> -49584.39

$ws.Cells.Item(137, 8).Value = 5182.7334   # H137: 5053.9355 -> 5182.7334
$ws.Cells.Item(137, 10).Value = 4834.3335   # J137: 4313.7144 -> 4834.3335
$ws.Cells.Item(137, 12).Value = 14503.0005   # L137: 12941.1432 -> 14503.0005
$ws.Cells.Item(137, 14).Value = -19603.0005   # N137: -18041.1432 -> -19603.0005

$ws.Cells.Item(138, 8).Value = 4682.101   # H138: 4637.45 -> 4682.101
$ws.Cells.Item(138, 9).Value = 3178.6428   # I138: 3039.6667 -> 3178.6428
$ws.Cells.Item(138, 10).Value = 5005.923   # J138: 5006.1694 -> 5005.923
$ws.Cells.Item(138, 11).Value = 9535.928400000001   # K138: 9119.000100000001 -> 9535.928400000001
$ws.Cells.Item(138, 12).Value = 15017.769   # L138: 15018.5082 -> 15017.769
$ws.Cells.Item(138, 13).Value = -4395.928400000001   # M138: -3979.000100000001 -> -4395.928400000001
$ws.Cells.Item(138, 14).Value = -25297.769   # N138: -25298.5082 -> -25297.769

$ws.Cells.Item(139, 8).Value = 50000   # H139: 49997.5 -> 50000
$ws.Cells.Item(139, 10).Value = 50000   # J139: 49997.5 -> 50000
$ws.Cells.Item(139, 12).Value = 50000   # L139: 49997.5 -> 50000
$ws.Cells.Item(139, 14).Value = -60280   # N139: -60277.5 -> -60280

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 3414.516   # H61: 3307.8125 -> 3414.516
$ws.Cells.Item(61, 9).Value = 3018.2693   # I61: 2908.0356 -> 3018.2693
$ws.Cells.Item(61, 10).Value = 5475   # J61: 6106.25 -> 5475
$ws.Cells.Item(61, 11).Value = 3018.2693   # K61: 2908.0356 -> 3018.2693
$ws.Cells.Item(61, 12).Value = 5475   # L61: 6106.25 -> 5475
$ws.Cells.Item(61, 13).Value = -2806.2693   # M61: -2696.0356 -> -2806.2693
$ws.Cells.Item(61, 14).Value = -5899   # N61: -6530.25 -> -5899

$ws.Cells.Item(74, 8).Value = 1609.2609   # H74: 1578.5385 -> 1609.2609
$ws.Cells.Item(74, 9).Value = 1527.2632   # I74: 1523.7826 -> 1527.2632
$ws.Cells.Item(74, 10).Value = 1998.75   # J74: 1998.3334 -> 1998.75
$ws.Cells.Item(74, 11).Value = 1527.2632   # K74: 1523.7826 -> 1527.2632
$ws.Cells.Item(74, 12).Value = 1998.75   # L74: 1998.3334 -> 1998.75
$ws.Cells.Item(74, 13).Value = -653.2632000000001   # M74: -649.7826 -> -653.2632000000001
$ws.Cells.Item(74, 14).Value = -3746.75   # N74: -3746.3334 -> -3746.75

$ws.Cells.Item(77, 8).Value = 1609.2609   # H77: 1578.5385 -> 1609.2609
$ws.Cells.Item(77, 9).Value = 1527.2632   # I77: 1523.7826 -> 1527.2632
$ws.Cells.Item(77, 10).Value = 1998.75   # J77: 1998.3334 -> 1998.75
$ws.Cells.Item(77, 11).Value = 7636.316000000001   # K77: 7618.913 -> 7636.316000000001
$ws.Cells.Item(77, 12).Value = 9993.75   # L77: 9991.666999999999 -> 9993.75
$ws.Cells.Item(77, 13).Value = -3268.316000000001   # M77: -3250.913 -> -3268.316000000001
$ws.Cells.Item(77, 14).Value = -18729.75   # N77: -18727.667 -> -18729.75

$ws.Cells.Item(109, 8).Value = 89999.664   # H109: 90000 -> 89999.664
$ws.Cells.Item(109, 10).Value = 89999.664   # J109: 90000 -> 89999.664
$ws.Cells.Item(109, 12).Value = 89999.664   # L109: 90000 -> 89999.664
$ws.Cells.Item(109, 14).Value = -92773.664   # N109: -92774 -> -92773.664

$ws.Cells.Item(122, 8).Value = 6204.8887   # H122: 6751.75 -> 6204.8887
$ws.Cells.Item(122, 9).Value = 7908.5   # I122: 11555.6 -> 7908.5
$ws.Cells.Item(122, 11).Value = 23725.5   # K122: 34666.8 -> 23725.5
$ws.Cells.Item(122, 13).Value = -21275.5   # M122: -32216.8 -> -21275.5

$ws.Cells.Item(132, 8).Value = 3580   # H132: 3421.6047 -> 3580
$ws.Cells.Item(132, 9).Value = 3285.5312   # I132: 3107.5881 -> 3285.5312
$ws.Cells.Item(132, 10).Value = 4627   # J132: 4607.8887 -> 4627
$ws.Cells.Item(132, 11).Value = 9856.5936   # K132: 9322.764299999999 -> 9856.5936
$ws.Cells.Item(132, 12).Value = 13881   # L132: 13823.6661 -> 13881
$ws.Cells.Item(132, 13).Value = -7326.5936   # M132: -6792.764299999999 -> -7326.5936
$ws.Cells.Item(132, 14).Value = -18941   # N132: -18883.6661 -> -18941

$ws.Cells.Item(136, 8).Value = 3414.516   # H136: 3307.8125 -> 3414.516
$ws.Cells.Item(136, 9).Value = 3018.2693   # I136: 2908.0356 -> 3018.2693
$ws.Cells.Item(136, 10).Value = 5475   # J136: 6106.25 -> 5475
$ws.Cells.Item(136, 11).Value = 9054.8079   # K136: 8724.106800000001 -> 9054.8079
$ws.Cells.Item(136, 12).Value = 16425   # L136: 18318.75 -> 16425
$ws.Cells.Item(136, 13).Value = -6504.8079   # M136: -6174.106800000001 -> -6504.8079
$ws.Cells.Item(136, 14).Value = -21525   # N136: -23418.75 -> -21525

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(99, 8).Value = 2435.7273   # H99: 2940.625 -> 2435.7273
$ws.Cells.Item(99, 9).Value = 2186.8572   # I99: 2807.8 -> 2186.8572
$ws.Cells.Item(99, 10).Value = 2871.25   # J99: 3162 -> 2871.25
$ws.Cells.Item(99, 11).Value = 2186.8572   # K99: 2807.8 -> 2186.8572
$ws.Cells.Item(99, 12).Value = 2871.25   # L99: 3162 -> 2871.25
$ws.Cells.Item(99, 13).Value = -688.8571999999999   # M99: -1309.8 -> -688.8571999999999
$ws.Cells.Item(99, 14).Value = -5867.25   # N99: -6158 -> -5867.25

$ws.Cells.Item(105, 8).Value = 112302   # H105: 167968 -> 112302
$ws.Cells.Item(105, 9).Value = 112302   # I105: 167968 -> 112302
$ws.Cells.Item(105, 11).Value = 112302   # K105: 167968 -> 112302
$ws.Cells.Item(105, 13).Value = -110555   # M105: -166221 -> -110555

$ws.Cells.Item(107, 8).Value = 717032   # H107: 836288 -> 717032
$ws.Cells.Item(107, 9).Value = 2105.4443   # I107: 2181.625 -> 2105.4443
$ws.Cells.Item(107, 10).Value = 2003899.8   # J107: 2504500.8 -> 2003899.8
$ws.Cells.Item(107, 11).Value = 2105.4443   # K107: 2181.625 -> 2105.4443
$ws.Cells.Item(107, 12).Value = 2003899.8   # L107: 2504500.8 -> 2003899.8
$ws.Cells.Item(107, 13).Value = -185.4443000000001   # M107: -261.625 -> -185.4443000000001
$ws.Cells.Item(107, 14).Value = -2007739.8   # N107: -2508340.8 -> -2007739.8

$ws.Cells.Item(134, 8).Value = 25876.043   # H134: 40379.965 -> 25876.043
$ws.Cells.Item(134, 9).Value = 4002.3572   # I134: 5952.8 -> 4002.3572
$ws.Cells.Item(134, 11).Value = 12007.0716   # K134: 17858.4 -> 12007.0716
$ws.Cells.Item(134, 13).Value = -9472.071599999999   # M134: -15323.4 -> -9472.071599999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 2884.9   # H58: 3138.7778 -> 2884.9
$ws.Cells.Item(58, 9).Value = 3745.8572   # I58: 4270.1665 -> 3745.8572
$ws.Cells.Item(58, 11).Value = 3745.8572   # K58: 4270.1665 -> 3745.8572
$ws.Cells.Item(58, 13).Value = -3542.8572   # M58: -4067.1665 -> -3542.8572

$ws.Cells.Item(136, 8).Value = 2884.9   # H136: 3138.7778 -> 2884.9
$ws.Cells.Item(136, 9).Value = 3745.8572   # I136: 4270.1665 -> 3745.8572
$ws.Cells.Item(136, 11).Value = 11237.5716   # K136: 12810.4995 -> 11237.5716
$ws.Cells.Item(136, 13).Value = -8687.571599999999   # M136: -10260.4995 -> -8687.571599999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 900920   # H5: 900928.3 -> 900920
$ws.Cells.Item(5, 9).Value = 81003.5   # I5: 81013.5 -> 81003.5
$ws.Cells.Item(5, 11).Value = 243010.5   # K5: 243040.5 -> 243010.5
$ws.Cells.Item(5, 13).Value = -242898.5   # M5: -242928.5 -> -242898.5

$ws.Cells.Item(80, 8).Value = 2000.1666   # H80: 1667 -> 2000.1666
$ws.Cells.Item(80, 10).Value = 2199.8   # J80: 1999.5 -> 2199.8
$ws.Cells.Item(80, 12).Value = 6599.400000000001   # L80: 5998.5 -> 6599.400000000001
$ws.Cells.Item(80, 14).Value = -8471.400000000001   # N80: -7870.5 -> -8471.400000000001

$ws.Cells.Item(83, 8).Value = 2000.1666   # H83: 1667 -> 2000.1666
$ws.Cells.Item(83, 10).Value = 2199.8   # J83: 1999.5 -> 2199.8
$ws.Cells.Item(83, 12).Value = 19798.2   # L83: 17995.5 -> 19798.2
$ws.Cells.Item(83, 14).Value = -29158.2   # N83: -27355.5 -> -29158.2

$ws.Cells.Item(107, 8).Value = 90043.74000000001   # H107: 86308.664 -> 90043.74000000001
$ws.Cells.Item(107, 10).Value = 147206.64   # J107: 137419.67 -> 147206.64
$ws.Cells.Item(107, 12).Value = 441619.92   # L107: 412259.01 -> 441619.92
$ws.Cells.Item(107, 14).Value = -445459.92   # N107: -416099.01 -> -445459.92

$ws.Cells.Item(116, 8).Value = 2579.6667   # H116: 3242.8333 -> 2579.6667
$ws.Cells.Item(116, 9).Value = 2354.5   # I116: 3091.4 -> 2354.5
$ws.Cells.Item(116, 10).Value = 3030   # J116: 4000 -> 3030
$ws.Cells.Item(116, 11).Value = 7063.5   # K116: 9274.200000000001 -> 7063.5
$ws.Cells.Item(116, 12).Value = 9090   # L116: 12000 -> 9090
$ws.Cells.Item(116, 13).Value = -3621.5   # M116: -5832.200000000001 -> -3621.5
$ws.Cells.Item(116, 14).Value = -15974   # N116: -18884 -> -15974

$ws.Cells.Item(127, 8).Value = 1338.6   # H127: 1483.7142 -> 1338.6
$ws.Cells.Item(127, 10).Value = 1338.6   # J127: 1483.7142 -> 1338.6
$ws.Cells.Item(127, 12).Value = 4015.8   # L127: 4451.142599999999 -> 4015.8
$ws.Cells.Item(127, 14).Value = -13935.8   # N127: -14371.1426 -> -13935.8

$ws.Cells.Item(131, 8).Value = 4554.5835   # H131: 4661.2085 -> 4554.5835
$ws.Cells.Item(131, 9).Value = 1462.6364   # I131: 1535.9 -> 1462.6364
$ws.Cells.Item(131, 10).Value = 7170.846   # J131: 6893.5713 -> 7170.846
$ws.Cells.Item(131, 11).Value = 4387.9092   # K131: 4607.700000000001 -> 4387.9092
$ws.Cells.Item(131, 12).Value = 21512.538   # L131: 20680.7139 -> 21512.538
$ws.Cells.Item(131, 13).Value = 652.0907999999999   # M131: 432.2999999999993 -> 652.0907999999999
$ws.Cells.Item(131, 14).Value = -31592.538   # N131: -30760.7139 -> -31592.538

$ws.Cells.Item(135, 8).Value = 900920   # H135: 900928.3 -> 900920
$ws.Cells.Item(135, 9).Value = 81003.5   # I135: 81013.5 -> 81003.5
$ws.Cells.Item(135, 11).Value = 729031.5   # K135: 729121.5 -> 729031.5
$ws.Cells.Item(135, 13).Value = -726496.5   # M135: -726586.5 -> -726496.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(20, 8).Value = 8826.362999999999   # H20: 9009 -> 8826.362999999999
$ws.Cells.Item(20, 9).Value = 7000   # I20: 0 -> 7000
$ws.Cells.Item(20, 11).Value = 7000   # K20: 0 -> 7000
$ws.Cells.Item(20, 13).Value = -6755   # M20: (none) -> -6755

$ws.Cells.Item(80, 8).Value = 1543639.8   # H80: 2005861.2 -> 1543639.8
$ws.Cells.Item(80, 9).Value = 1004490.5   # I80: 1254888.2 -> 1004490.5
$ws.Cells.Item(80, 10).Value = 3340803.8   # J80: 5009753 -> 3340803.8
$ws.Cells.Item(80, 11).Value = 1004490.5   # K80: 1254888.2 -> 1004490.5
$ws.Cells.Item(80, 12).Value = 3340803.8   # L80: 5009753 -> 3340803.8
$ws.Cells.Item(80, 13).Value = -1003492.5   # M80: -1253890.2 -> -1003492.5
$ws.Cells.Item(80, 14).Value = -3342799.8   # N80: -5011749 -> -3342799.8

$ws.Cells.Item(83, 8).Value = 1543639.8   # H83: 2005861.2 -> 1543639.8
$ws.Cells.Item(83, 9).Value = 1004490.5   # I83: 1254888.2 -> 1004490.5
$ws.Cells.Item(83, 10).Value = 3340803.8   # J83: 5009753 -> 3340803.8
$ws.Cells.Item(83, 11).Value = 5022452.5   # K83: 6274441 -> 5022452.5
$ws.Cells.Item(83, 12).Value = 16704019   # L83: 25048765 -> 16704019
$ws.Cells.Item(83, 13).Value = -5017460.5   # M83: -6269449 -> -5017460.5
$ws.Cells.Item(83, 14).Value = -16714003   # N83: -25058749 -> -16714003

$ws.Cells.Item(93, 8).Value = 39957.4   # H93: 39956.332 -> 39957.4
$ws.Cells.Item(93, 10).Value = 39959.5   # J93: 39960 -> 39959.5
$ws.Cells.Item(93, 12).Value = 39959.5   # L93: 39960 -> 39959.5
$ws.Cells.Item(93, 14).Value = -43703.5   # N93: -43704 -> -43703.5

$ws.Cells.Item(123, 8).Value = 42499.25   # H123: 38999.4 -> 42499.25
$ws.Cells.Item(123, 10).Value = 42499.25   # J123: 38999.4 -> 42499.25
$ws.Cells.Item(123, 12).Value = 42499.25   # L123: 38999.4 -> 42499.25
$ws.Cells.Item(123, 14).Value = -47399.25   # N123: -43899.4 -> -47399.25

$ws.Cells.Item(132, 8).Value = 42179.82   # H132: 38444.87 -> 42179.82
$ws.Cells.Item(132, 9).Value = 6001.346   # I132: 6616.1665 -> 6001.346
$ws.Cells.Item(132, 10).Value = 512500   # J132: 147571.86 -> 512500
$ws.Cells.Item(132, 11).Value = 18004.038   # K132: 19848.4995 -> 18004.038
$ws.Cells.Item(132, 12).Value = 1537500   # L132: 442715.58 -> 1537500
$ws.Cells.Item(132, 13).Value = -15474.038   # M132: -17318.4995 -> -15474.038
$ws.Cells.Item(132, 14).Value = -1542560   # N132: -447775.58 -> -1542560

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 8798.143   # H7: 6407.727 -> 8798.143
$ws.Cells.Item(7, 9).Value = 9473   # I7: 5848.75 -> 9473
$ws.Cells.Item(7, 11).Value = 9473   # K7: 5848.75 -> 9473
$ws.Cells.Item(7, 13).Value = -9361   # M7: -5736.75 -> -9361

$ws.Cells.Item(23, 8).Value = 641666.7   # H23: 578428.5600000001 -> 641666.7
$ws.Cells.Item(23, 9).Value = 641666.7   # I23: 578428.5600000001 -> 641666.7
$ws.Cells.Item(23, 11).Value = 641666.7   # K23: 578428.5600000001 -> 641666.7
$ws.Cells.Item(23, 13).Value = -641436.7   # M23: -578198.5600000001 -> -641436.7

$ws.Cells.Item(61, 8).Value = 5637.722   # H61: 6070.5 -> 5637.722
$ws.Cells.Item(61, 9).Value = 6505.3   # I61: 6984.3335 -> 6505.3
$ws.Cells.Item(61, 10).Value = 4553.25   # J61: 4895.5713 -> 4553.25
$ws.Cells.Item(61, 11).Value = 6505.3   # K61: 6984.3335 -> 6505.3
$ws.Cells.Item(61, 12).Value = 4553.25   # L61: 4895.5713 -> 4553.25
$ws.Cells.Item(61, 13).Value = -6303.3   # M61: -6782.3335 -> -6303.3
$ws.Cells.Item(61, 14).Value = -4957.25   # N61: -5299.5713 -> -4957.25

$ws.Cells.Item(113, 8).Value = 5637.722   # H113: 6070.5 -> 5637.722
$ws.Cells.Item(113, 9).Value = 6505.3   # I113: 6984.3335 -> 6505.3
$ws.Cells.Item(113, 10).Value = 4553.25   # J113: 4895.5713 -> 4553.25
$ws.Cells.Item(113, 11).Value = 6505.3   # K113: 6984.3335 -> 6505.3
$ws.Cells.Item(113, 12).Value = 4553.25   # L113: 4895.5713 -> 4553.25
$ws.Cells.Item(113, 13).Value = -4335.3   # M113: -4814.3335 -> -4335.3
$ws.Cells.Item(113, 14).Value = -8893.25   # N113: -9235.5713 -> -8893.25

$ws.Cells.Item(122, 8).Value = 5011.353   # H122: 5439.533 -> 5011.353
$ws.Cells.Item(122, 10).Value = 5356.4287   # J122: 6779 -> 5356.4287
$ws.Cells.Item(122, 12).Value = 16069.2861   # L122: 20337 -> 16069.2861
$ws.Cells.Item(122, 14).Value = -20969.2861   # N122: -25237 -> -20969.2861

$ws.Cells.Item(126, 8).Value = 8798.143   # H126: 6407.727 -> 8798.143
$ws.Cells.Item(126, 9).Value = 9473   # I126: 5848.75 -> 9473
$ws.Cells.Item(126, 11).Value = 28419   # K126: 17546.25 -> 28419
$ws.Cells.Item(126, 13).Value = -25949   # M126: -15076.25 -> -25949

$ws.Cells.Item(136, 8).Value = 377652.44   # H136: 392038.88 -> 377652.44
$ws.Cells.Item(136, 10).Value = 9334.267   # J136: 9743.5 -> 9334.267
$ws.Cells.Item(136, 12).Value = 28002.801   # L136: 29230.5 -> 28002.801
$ws.Cells.Item(136, 14).Value = -33102.801   # N136: -34330.5 -> -33102.801

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(9, 8).Value = 0   # H9: 99999 -> 0
$ws.Cells.Item(9, 9).Value = 0   # I9: 99999 -> 0
$ws.Cells.Item(9, 11).Value = 0   # K9: 99999 -> 0
$ws.Cells.Item(9, 13).ClearContents()   # M9: removed (no NQ profit / placeholder row)

$ws.Cells.Item(37, 8).Value = 0   # H37: 1000 -> 0
$ws.Cells.Item(37, 9).Value = 0   # I37: 1000 -> 0
$ws.Cells.Item(37, 11).Value = 0   # K37: 1000 -> 0
$ws.Cells.Item(37, 13).ClearContents()   # M37: removed (no NQ profit / placeholder row)

$ws.Cells.Item(46, 8).Value = 75887.25   # H46: 79999.664 -> 75887.25
$ws.Cells.Item(46, 10).Value = 75887.25   # J46: 79999.664 -> 75887.25
$ws.Cells.Item(46, 12).Value = 75887.25   # L46: 79999.664 -> 75887.25
$ws.Cells.Item(46, 14).Value = -76349.25   # N46: -80461.664 -> -76349.25

$ws.Cells.Item(107, 8).Value = 977.3333   # H107: 775.7692 -> 977.3333
$ws.Cells.Item(107, 9).Value = 1002.0625   # I107: 818.4545000000001 -> 1002.0625
$ws.Cells.Item(107, 10).Value = 779.5   # J107: 541 -> 779.5
$ws.Cells.Item(107, 11).Value = 3006.1875   # K107: 2455.3635 -> 3006.1875
$ws.Cells.Item(107, 12).Value = 2338.5   # L107: 1623 -> 2338.5
$ws.Cells.Item(107, 13).Value = -1086.1875   # M107: -535.3635000000004 -> -1086.1875
$ws.Cells.Item(107, 14).Value = -6178.5   # N107: -5463 -> -6178.5

$ws.Cells.Item(122, 8).Value = 58827024   # H122: 55558970 -> 58827024
$ws.Cells.Item(122, 10).Value = 5340.2   # J122: 4784.3335 -> 5340.2
$ws.Cells.Item(122, 12).Value = 16020.6   # L122: 14353.0005 -> 16020.6
$ws.Cells.Item(122, 14).Value = -20920.6   # N122: -19253.0005 -> -20920.6

$ws.Cells.Item(132, 8).Value = 26960.256   # H132: 28975.625 -> 26960.256
$ws.Cells.Item(132, 9).Value = 1363.7826   # I132: 1555.05 -> 1363.7826
$ws.Cells.Item(132, 11).Value = 4091.3478   # K132: 4665.15 -> 4091.3478
$ws.Cells.Item(132, 13).Value = -1561.3478   # M132: -2135.15 -> -1561.3478

$ws.Cells.Item(134, 8).Value = 75887.25   # H134: 79999.664 -> 75887.25
$ws.Cells.Item(134, 10).Value = 75887.25   # J134: 79999.664 -> 75887.25
$ws.Cells.Item(134, 12).Value = 227661.75   # L134: 239998.992 -> 227661.75
$ws.Cells.Item(134, 14).Value = -232731.75   # N134: -245068.992 -> -232731.75

$ws.Cells.Item(136, 8).Value = 20962392   # H136: 22359752 -> 20962392
$ws.Cells.Item(136, 9).Value = 25645482   # I136: 27782440 -> 25645482
$ws.Cells.Item(136, 11).Value = 76936446   # K136: 83347320 -> 76936446
$ws.Cells.Item(136, 13).Value = -76933896   # M136: -83344770 -> -76933896
